# "Tried to implement Penality Reward System (unfinished)"
#
# Weekly Quantity sheet: two weekly buckets (the weeks of 45123.99999999999
# and 45137.99999999999) are dropped from the "Order Week" series, which
# shifts every later row up by two and shrinks the used range from
# A1:B47 to A1:B45. The week that stays at row 21 gets its requested
# quantity corrected from 36 to 34, and the row that lands on 22 after the
# shift (week 45130.99999999999) gets its quantity corrected to 64.
#
# Monthly Trend sheet: the month total for 45138.99999999999 is corrected
# from 192 down to 98.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Weekly Quantity")

# Delete the higher row index first so the second Delete() still targets
# the intended row (row 22 here refers to the original sheet row, which
# after deleting row 24 is still row 22 - the "45123.99999999999" week).
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(22).Delete()

# Fix up the requested-quantity values that changed (not just shifted).
$ws.Range("B21").Value = 34
$ws.Range("B22").Value = 64

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B8").Value = 98
